$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the "Macroferia Regional de Talca -
# Espárragos" sheet, inserted above the previous most-recent entry (old
# row 55). Inserting a whole row there shifts every following row down by
# one (old 55-62 -> new 56-63) while preserving their contents untouched,
# and the sheet's used range grows from A1:R62 to A1:R63 automatically -
# exactly matching the target diff.
$ws.Rows(55).Insert()

# Populate the newly inserted row 55 with the new weekly entry.
$ws.Range("A55").Value = 5
$ws.Range("B55").Value = "Macroferia Regional de Talca"
$ws.Range("C55").Value = "Maule"
$ws.Range("D55").Value = 44826
$ws.Range("E55").Value = 7
$ws.Range("F55").Value = 300000000
$ws.Range("G55").Value = "Espárragos"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 2000
$ws.Range("M55").Value = 2000
$ws.Range("N55").Value = "$/kilo"
$ws.Range("O55").Value = "Provincia de Linares"
$ws.Range("P55").Value = 2000
$ws.Range("Q55").Value = 1
$ws.Range("R55").Value = "Hortaliza"
